$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row, day, month, year, pruebas realizadas, pruebas positivas)
$rows = @(
    @(277, 29, 12, 2020, 5403, 879),
    @(278, 30, 12, 2020, 4613, 846),
    @(279, 31, 12, 2020, 1115, 224),
    @(280, 1, 1, 2021, 685, 80),
    @(281, 2, 1, 2021, 656, 159),
    @(282, 3, 1, 2021, 631, 181),
    @(283, 4, 1, 2021, 5132, 763),
    @(284, 5, 1, 2021, 4741, 783),
    @(285, 6, 1, 2021, 5637, 872),
    @(286, 7, 1, 2021, 5583, 990),
    @(287, 8, 1, 2021, 5998, 1063),
    @(288, 9, 1, 2021, 1747, 116),
    @(289, 10, 1, 2021, 883, 172)
)

foreach ($r in $rows) {
    $row = $r[0]
    $day = $r[1]
    $month = $r[2]
    $year = $r[3]
    $pruebas = $r[4]
    $positivas = $r[5]

    # Copy formatting from the row directly above (keeps shared style indices)
    $ws.Range("A$($row-1):F$($row-1)").Copy()
    $ws.Range("A$($row):F$($row)").PasteSpecial(-4122)
    $ws.Application.CutCopyMode = 0
    $ws.Rows.Item($row).RowHeight = 15

    $ws.Cells.Item($row, 2).Value = $day
    $ws.Cells.Item($row, 3).Value = $month
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = $pruebas
    $ws.Cells.Item($row, 6).Value = $positivas

    $f = '=+Hoja1!$B' + $row + '&"/"&Hoja1!$C' + $row + '&"/"&Hoja1!$D' + $row
    $ws.Cells.Item($row, 1).Formula = $f
}

# Expand the table (ListObject) to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I289"))

# Update the visible selection to match the post-edit state
$null = $excel.Goto($ws.Range("A265"), $true)
$null = $ws.Range("F290").Select()
